$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ChooseFormula($row) {
    $f = '=CHOOSE(D' + $row + ',"Chapter DLC","Half-Chapter DLC","Clothing Pack DLC","Original Soundtrack DLC","Character Pack DLC","Other","Retracted","Chapter Pack DLC")'
    $ws.Range("E" + $row).Formula = $f
}

# Row 44 already exists (id 56 / Tokyo Ghoul / Icecream) - fill in the
# previously-missing cells: rDate, category, dlcCategories formula, release
$ws.Range("C44").Value = "'02.04.2025"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = 1
Set-ChooseFormula 44
$ws.Range("F44").Value = "8.6.0"

# Row 45 - brand new DLC row (id 57, Jerky / Steady Pulse / Orela Rose)
$ws.Range("A45").Value = 57
$ws.Range("B45").Value = "Steady Pulse"
$ws.Range("C45").Value = "'06.05.2025"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = 2
Set-ChooseFormula 45
$ws.Range("F45").Value = "8.7.0"
$ws.Range("G45").Value = "Jerky"
$ws.Range("H45").Value = "Jerky"
$ws.Range("I45").Value = "Orela Rose"

# Row 46 - brand new DLC row (id 58, Five Nights at Freddy's / TBA)
$ws.Range("A46").Value = 58
$ws.Range("B46").Value = "Five Nights at Freddy's"
$ws.Range("C46").Value = "'17.06.2025"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = 1
Set-ChooseFormula 46
$ws.Range("F46").Value = "9.0.0"
$ws.Range("G46").Value = "TBA"
$ws.Range("H46").Value = "TBA"
$ws.Range("I46").Value = "TBA"

# Match the author's final on-screen selection / scroll position
$ws.Range("A16").Select() | Out-Null
$ws.Range("E46").Select() | Out-Null

Write-Host "done"
